$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 497.07144
$ws.Range("I9").Value = 180.22223
$ws.Range("J9").Value = 1067.4
$ws.Range("K9").Value = 180.22223
$ws.Range("L9").Value = 1067.4
$ws.Range("M9").Value = -11.22223
$ws.Range("N9").Value = -1405.4

$ws.Range("H15").Value = 935.6316
$ws.Range("I15").Value = 935.6316
$ws.Range("K15").Value = 2806.8948
$ws.Range("M15").Value = -2637.8948

$ws.Range("H32").Value = 5000
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 5000
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 5000
$ws.Range("M32").Value = ""
$ws.Range("N32").Value = -5652

$ws.Range("H86").Value = 5877.1
$ws.Range("I86").Value = 2700
$ws.Range("J86").Value = 7238.7144
$ws.Range("K86").Value = 2700
$ws.Range("L86").Value = 7238.7144
$ws.Range("M86").Value = -1577
$ws.Range("N86").Value = -9484.714400000001

$ws.Range("H89").Value = 5877.1
$ws.Range("I89").Value = 2700
$ws.Range("J89").Value = 7238.7144
$ws.Range("K89").Value = 13500
$ws.Range("L89").Value = 36193.572
$ws.Range("M89").Value = -7884
$ws.Range("N89").Value = -47425.572

$ws.Range("H93").Value = 30000
$ws.Range("J93").Value = 30000
$ws.Range("L93").Value = 30000
$ws.Range("N93").Value = -34992

$ws.Range("H141").Value = 1857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1156.6875
$ws.Range("I2").Value = 607.6429000000001
$ws.Range("K2").Value = 607.6429000000001
$ws.Range("M2").Value = -494.6429000000001

$ws.Range("H32").Value = 5144.879
$ws.Range("I32").Value = 4702.645
$ws.Range("K32").Value = 4702.645
$ws.Range("M32").Value = -4415.645

$ws.Range("H61").Value = 3254.6924
$ws.Range("I61").Value = 1506.7894
$ws.Range("K61").Value = 1506.7894
$ws.Range("M61").Value = -1294.7894

$ws.Range("H102").Value = 892.4666999999999
$ws.Range("I102").Value = 907.7692
$ws.Range("K102").Value = 907.7692
$ws.Range("M102").Value = 714.2308

$ws.Range("H116").Value = 1156.6875
$ws.Range("I116").Value = 607.6429000000001
$ws.Range("K116").Value = 607.6429000000001
$ws.Range("M116").Value = 1686.3571

$ws.Range("H117").Value = 49999
$ws.Range("J117").Value = 49999
$ws.Range("L117").Value = 49999
$ws.Range("N117").Value = -59177

$ws.Range("H122").Value = 2161
$ws.Range("I122").Value = 1707
$ws.Range("J122").Value = 2766.3333
$ws.Range("K122").Value = 5121
$ws.Range("L122").Value = 8298.999899999999
$ws.Range("M122").Value = -2671
$ws.Range("N122").Value = -13198.9999

$ws.Range("H136").Value = 3254.6924
$ws.Range("I136").Value = 1506.7894
$ws.Range("K136").Value = 4520.3682
$ws.Range("M136").Value = -1970.3682

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1156.6875
$ws.Range("I3").Value = 607.6429000000001
$ws.Range("K3").Value = 607.6429000000001
$ws.Range("M3").Value = -493.6429000000001

$ws.Range("H107").Value = 4739.1665
$ws.Range("I107").Value = 3102.25
$ws.Range("K107").Value = 3102.25
$ws.Range("M107").Value = -1182.25

$ws.Range("H135").Value = 99995
$ws.Range("J135").Value = 99995
$ws.Range("L135").Value = 99995
$ws.Range("N135").Value = -110135

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 1450
$ws.Range("I38").Value = 1450
$ws.Range("K38").Value = 1450
$ws.Range("M38").Value = -1073

$ws.Range("H46").Value = 1450
$ws.Range("I46").Value = 1450
$ws.Range("K46").Value = 1450
$ws.Range("M46").Value = -1239

$ws.Range("H58").Value = 2213.125
$ws.Range("J58").Value = 2577.6667
$ws.Range("L58").Value = 2577.6667
$ws.Range("N58").Value = -2983.6667

$ws.Range("H88").Value = 12750.777
$ws.Range("J88").Value = 12750.777
$ws.Range("L88").Value = 12750.777
$ws.Range("N88").Value = -13562.777

$ws.Range("H91").Value = 12750.777
$ws.Range("J91").Value = 12750.777
$ws.Range("L91").Value = 12750.777
$ws.Range("N91").Value = -15558.777

$ws.Range("H105").Value = 3496.375
$ws.Range("I105").Value = 2723.3845
$ws.Range("J105").Value = 4409.909
$ws.Range("K105").Value = 2723.3845
$ws.Range("L105").Value = 4409.909
$ws.Range("M105").Value = -976.3845000000001
$ws.Range("N105").Value = -7903.909

$ws.Range("H132").Value = 1356.6
$ws.Range("I132").Value = 1140.125
$ws.Range("K132").Value = 3420.375
$ws.Range("M132").Value = -890.375

$ws.Range("H136").Value = 2213.125
$ws.Range("J136").Value = 2577.6667
$ws.Range("L136").Value = 7733.000100000001
$ws.Range("N136").Value = -12833.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 4596.6
$ws.Range("I114").Value = 10000
$ws.Range("J114").Value = 3245.75
$ws.Range("K114").Value = 30000
$ws.Range("L114").Value = 9737.25
$ws.Range("M114").Value = -26746
$ws.Range("N114").Value = -16245.25

$ws.Range("H131").Value = 2939.7646
$ws.Range("I131").Value = 1997.5
$ws.Range("K131").Value = 5992.5
$ws.Range("M131").Value = -952.5

$ws.Range("H137").Value = 4259.222
$ws.Range("J137").Value = 5374
$ws.Range("L137").Value = 16122
$ws.Range("N137").Value = -26322

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").Value = ""

$ws.Range("H33").Value = 19999.5
$ws.Range("J33").Value = 19999.5
$ws.Range("L33").Value = 19999.5
$ws.Range("N33").Value = -20503.5

$ws.Range("H122").Value = 2303.0588
$ws.Range("I122").Value = 1237.1666
$ws.Range("K122").Value = 3711.4998
$ws.Range("M122").Value = -1261.4998

$ws.Range("H132").Value = 1759
$ws.Range("I132").Value = 1961.6666
$ws.Range("K132").Value = 5884.9998
$ws.Range("M132").Value = -3354.9998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2763
$ws.Range("I68").Value = 2095.25
$ws.Range("J68").Value = 3144.5715
$ws.Range("K68").Value = 2095.25
$ws.Range("L68").Value = 3144.5715
$ws.Range("M68").Value = -1346.25
$ws.Range("N68").Value = -4642.5715

$ws.Range("H71").Value = 2763
$ws.Range("I71").Value = 2095.25
$ws.Range("J71").Value = 3144.5715
$ws.Range("K71").Value = 10476.25
$ws.Range("L71").Value = 15722.8575
$ws.Range("M71").Value = -6732.25
$ws.Range("N71").Value = -23210.8575

$ws.Range("H136").Value = 1483.3334
$ws.Range("I136").Value = 1275
$ws.Range("J136").Value = 1900
$ws.Range("K136").Value = 3825
$ws.Range("L136").Value = 5700
$ws.Range("M136").Value = -1275
$ws.Range("N136").Value = -10800

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 81531
$ws.Range("J46").Value = 81531
$ws.Range("L46").Value = 81531
$ws.Range("N46").Value = -81993

$ws.Range("H134").Value = 81531
$ws.Range("J134").Value = 81531
$ws.Range("L134").Value = 244593
$ws.Range("N134").Value = -249663

$ws.Range("H136").Value = 2204.3809
$ws.Range("J136").Value = 1514
$ws.Range("L136").Value = 4542
$ws.Range("N136").Value = -9642
